# This edit rotates the data contained in the 5 data rows (rows 2-6,
# directly below the header row) of the active sheet. Each row's full
# contents (columns A through AY) move down to a new row as follows:
#
#   old row 5 -> new row 2
#   old row 6 -> new row 3
#   old row 2 -> new row 4
#   old row 3 -> new row 5
#   old row 4 -> new row 6
#
# i.e. a cyclic rotation of the 5 records. No cell values themselves are
# changed - only their row position moves.
#
# Approach: snapshot every source row into a scratch area further down
# the sheet using Range.Copy (which -- unlike assigning .Value() --
# preserves the original cell typing, e.g. numeric-looking text such as
# "1" or date-looking text such as "2020-05-22" stay text instead of
# being re-interpreted by Excel). The destination rows are cleared
# before pasting into them, because copying a blank source cell over an
# already-populated destination cell does not blank the destination.
# Finally the scratch rows are removed again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1    # column A
$lastCol  = 51   # column AY (last used column on the sheet)

$firstDataRow = 2
$lastDataRow  = 6

# destination row -> source row
$rowMap = @{
    2 = 5
    3 = 6
    4 = 2
    5 = 3
    6 = 4
}

# Copy each data row down into a scratch area (rows far below the data)
# so the rotation can be performed without a source row being
# overwritten before it has been read.
$scratchBase = 1000
$scratchRowOf = @{}
$offset = 0
foreach ($srcRow in $firstDataRow..$lastDataRow) {
    $offset++
    $scratchRow = $scratchBase + $offset
    $scratchRowOf[$srcRow] = $scratchRow

    $srcRange = $ws.Range($ws.Cells.Item($srcRow, $firstCol), $ws.Cells.Item($srcRow, $lastCol))
    $destCell = $ws.Cells.Item($scratchRow, $firstCol)
    $srcRange.Copy($destCell)
}

# Write each scratch row into its final destination row.
foreach ($destRow in $firstDataRow..$lastDataRow) {
    $srcRow = $rowMap[$destRow]
    $scratchRow = $scratchRowOf[$srcRow]

    $destRange = $ws.Range($ws.Cells.Item($destRow, $firstCol), $ws.Cells.Item($destRow, $lastCol))
    $destRange.ClearContents()

    $scratchRange = $ws.Range($ws.Cells.Item($scratchRow, $firstCol), $ws.Cells.Item($scratchRow, $lastCol))
    $destCell = $ws.Cells.Item($destRow, $firstCol)
    $scratchRange.Copy($destCell)
}

# Remove the temporary scratch rows again.
$scratchRange = $ws.Range($ws.Cells.Item($scratchBase + 1, $firstCol), $ws.Cells.Item($scratchBase + 5, $lastCol))
$scratchRange.ClearContents()
